$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (A1:L1) - Room/Block/Group reordered, Level & Course added
$header = New-Object 'object[,]' 1,12
$header[0,0]  = "Day"
$header[0,1]  = "Time"
$header[0,2]  = "Module Code"
$header[0,3]  = "Module Title"
$header[0,4]  = "Hours"
$header[0,5]  = "Class Type"
$header[0,6]  = "Lecturer"
$header[0,7]  = "Room"
$header[0,8]  = "Block"
$header[0,9]  = "Group"
$header[0,10] = "Level"
$header[0,11] = "Course"
$ws.Range("A1:L1").Value = $header

# Full replacement data for rows 2-10 (A:L)
$data = @(
    @("SUN", "9:30-11:30",  "5CS024", "Collaborative Development",                 2,   "Tutorial", "Mr. Udaya Kandel",   "SR-02 Bilston",  "WLV", "L5CG6",         5, "BCS"),
    @("SUN", "12:30-15:00", "5CS022", "Human Computer Interaction",                2.5, "Workshop", "Mr. Dipesh Shrestha","SR-04 Crompton", "WLV", "L5CG6",         5, "BCS"),
    @("MON", "9:30-12:00",  "5CS024", "Collaborative Development",                 2.5, "Workshop", "Mr. Udaya Kandel",   "SR-04 Crompton", "WLV", "L5CG6",         5, "BCS"),
    @("TUE", "7:00-9:00",   "5CS022", "Human Computer Interaction",                2,   "Lecture",  "Mr. Apurba Neupane", "LT-02 Telford",  "WLV", "L5CG(5+6+7+8)", 5, "BCS"),
    @("TUE", "9:30-11:30",  "5CS020", "Distributed and Cloud Systems Programming", 2,   "Lecture",  "Mr. Sumanta Silwal", "LT-01 Wulfruna", "WLV", "L5CG(5+6+7+8)", 5, "BCS"),
    @("WED", "7:00-9:00",   "5CS024", "Collaborative Development",                 2,   "Lecture",  "Mr. Raj Shrestha",   "LT-02 Telford",  "WLV", "L5CG(5+6+7+8)", 5, "BCS"),
    @("WED", "9:30-11:30",  "5CS022", "Human Computer Interaction",                2,   "Tutorial", "Mr. Dipesh Shrestha","SR-02 Bilston",  "WLV", "L5CG6",         5, "BCS"),
    @("THU", "13:00-15:00", "5CS020", "Distributed and Cloud Systems Programming", 2,   "Tutorial", "Mr. Prabin Sapkota", "SR-02 Bilston",  "WLV", "L5CG6",         5, "BCS"),
    @("FRI", "12:30-15:30", "5CS020", "Distributed and Cloud Systems Programming", 2.5, "Workshop", "Mr. Prabin Sapkota", "Lab-01 Mander",  "WLV", "L5CG6",         5, "BCS")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowArr = New-Object 'object[,]' 1,12
    for ($c = 0; $c -lt 12; $c++) {
        $rowArr[0,$c] = $data[$i][$c]
    }
    $ws.Range("A${row}:L${row}").Value = $rowArr
}
